# Update the reporting period (4to trimestre 2020) for rows 9-60 on the
# main "Reporte de Formatos" sheet.
#   B column: "Fecha de inicio del periodo que se informa" 44013 -> 44105
#   C column: "Fecha de termino del periodo que se informa" 44104 -> 44196
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reporte de Formatos")

for ($r = 9; $r -le 60; $r++) {
    $ws.Cells.Item($r, 2).Value2 = 44105
    $ws.Cells.Item($r, 3).Value2 = 44196
}

# Move the visible selection / scroll position that Excel persisted the
# last time the workbook was saved.
$ws.Activate()
$ws.Range("D65").Select()
